$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new cell value - this also fixes the dimension and shared strings automatically
$ws.Range("E11").Value = "s"

# Update the active selection to the newly edited cell
$ws.Range("E11").Select()
